$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registros")

# Update the street address in I2 (added "B" after the house number 1881)
$ws.Range("I2").Value = "Rua Tereza Feitosa, 1881B, Nossa Senhora das Gracas, 64519-410, Teresina-PI"

# Narrow most of the columns (columns D and H are left untouched)
$ws.Columns.Item(1).ColumnWidth = 4.666666666666667
$ws.Columns.Item(2).ColumnWidth = 7.833333333333333
$ws.Columns.Item(3).ColumnWidth = 5.5
$ws.Columns.Item(5).ColumnWidth = 5.666666666666667
$ws.Columns.Item(6).ColumnWidth = 6.666666666666667
$ws.Columns.Item(7).ColumnWidth = 14.666666666666666
$ws.Columns.Item(9).ColumnWidth = 8.333333333333334
$ws.Columns.Item(10).ColumnWidth = 7.5

# Move the active selection from J5 to N5
$ws.Range("N5").Select()
